# Refresh market-data-driven columns (H-N) across the Leve Profit tables.
# Generated from the authoritative before/after cell diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 937.5
$ws.Range("I43").Value = 923.5
$ws.Range("J43").Value = 951.5
$ws.Range("K43").Value = 923.5
$ws.Range("L43").Value = 951.5
$ws.Range("M43").Value = -854.5
$ws.Range("N43").Value = -1089.5

$ws.Range("H70").Value = 1796.4814
$ws.Range("I70").Value = 1621.0526
$ws.Range("J70").Value = 2213.125
$ws.Range("K70").Value = 4863.1578
$ws.Range("L70").Value = 6639.375
$ws.Range("M70").Value = -4593.1578
$ws.Range("N70").Value = -7179.375

$ws.Range("H73").Value = 1796.4814
$ws.Range("I73").Value = 1621.0526
$ws.Range("J73").Value = 2213.125
$ws.Range("K73").Value = 4863.1578
$ws.Range("L73").Value = 6639.375
$ws.Range("M73").Value = -3927.1578
$ws.Range("N73").Value = -8511.375

$ws.Range("H116").Value = 2179.1936
$ws.Range("I116").Value = 1950.5
$ws.Range("J116").Value = 2495.8462
$ws.Range("K116").Value = 1950.5
$ws.Range("L116").Value = 2495.8462
$ws.Range("M116").Value = 1491.5
$ws.Range("N116").Value = -9379.8462

$ws.Range("H138").Value = 3635.205
$ws.Range("I138").Value = 974.75
$ws.Range("J138").Value = 4817.6294
$ws.Range("K138").Value = 2924.25
$ws.Range("L138").Value = 14452.8882
$ws.Range("M138").Value = 2215.75
$ws.Range("N138").Value = -24732.8882

$ws.Range("H141").Value = 4411.4287
$ws.Range("I141").Value = 4470
$ws.Range("J141").Value = 4333.3335
$ws.Range("K141").Value = 13410
$ws.Range("L141").Value = 13000.0005
$ws.Range("M141").Value = -8230
$ws.Range("N141").Value = -23360.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1569.674
$ws.Range("I45").Value = 1499.7273
$ws.Range("J45").Value = 1747.2307
$ws.Range("K45").Value = 1499.7273
$ws.Range("L45").Value = 1747.2307
$ws.Range("M45").Value = -1122.7273
$ws.Range("N45").Value = -2501.2307

$ws.Range("H76").Value = 36890.555
$ws.Range("I76").Value = 10000
$ws.Range("J76").Value = 40251.875
$ws.Range("K76").Value = 10000
$ws.Range("L76").Value = 40251.875
$ws.Range("M76").Value = -9662
$ws.Range("N76").Value = -40927.875

$ws.Range("H79").Value = 36890.555
$ws.Range("I79").Value = 10000
$ws.Range("J79").Value = 40251.875
$ws.Range("K79").Value = 10000
$ws.Range("L79").Value = 40251.875
$ws.Range("M79").Value = -8830
$ws.Range("N79").Value = -42591.875

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").ClearContents()
$ws.Range("N114").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1971.1964
$ws.Range("I86").Value = 1832.3914
$ws.Range("J86").Value = 2609.7
$ws.Range("K86").Value = 1832.3914
$ws.Range("L86").Value = 2609.7
$ws.Range("M86").Value = -709.3914
$ws.Range("N86").Value = -4855.7

$ws.Range("H89").Value = 1971.1964
$ws.Range("I89").Value = 1832.3914
$ws.Range("J89").Value = 2609.7
$ws.Range("K89").Value = 9161.957
$ws.Range("L89").Value = 13048.5
$ws.Range("M89").Value = -3545.957
$ws.Range("N89").Value = -24280.5

$ws.Range("H134").Value = 8112.5
$ws.Range("I134").Value = 7482.1665
$ws.Range("K134").Value = 22446.4995
$ws.Range("M134").Value = -19911.4995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 66.666664
$ws.Range("I7").Value = 66.666664
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 66.666664
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = 46.333336

$ws.Range("H31").Value = 6794.9824
$ws.Range("I31").Value = 5070.125
$ws.Range("J31").Value = 8049.4243
$ws.Range("K31").Value = 5070.125
$ws.Range("L31").Value = 8049.4243
$ws.Range("M31").Value = -4775.125
$ws.Range("N31").Value = -8639.424299999999

$ws.Range("H34").Value = 6794.9824
$ws.Range("I34").Value = 5070.125
$ws.Range("J34").Value = 8049.4243
$ws.Range("K34").Value = 5070.125
$ws.Range("L34").Value = 8049.4243
$ws.Range("M34").Value = -4868.125
$ws.Range("N34").Value = -8453.424299999999

$ws.Range("H59").Value = 23040.54
$ws.Range("J59").Value = 23040.54
$ws.Range("L59").Value = 23040.54
$ws.Range("N59").Value = -25330.54

$ws.Range("H68").Value = 20295
$ws.Range("J68").Value = 20295
$ws.Range("L68").Value = 20295
$ws.Range("N68").Value = -21793

$ws.Range("H71").Value = 20295
$ws.Range("J71").Value = 20295
$ws.Range("L71").Value = 60885
$ws.Range("N71").Value = -68373

$ws.Range("H74").Value = 28045.4
$ws.Range("I74").Value = 11642.5
$ws.Range("J74").Value = 38980.668
$ws.Range("K74").Value = 11642.5
$ws.Range("L74").Value = 38980.668
$ws.Range("M74").Value = -10768.5
$ws.Range("N74").Value = -40728.668

$ws.Range("H77").Value = 28045.4
$ws.Range("I77").Value = 11642.5
$ws.Range("J77").Value = 38980.668
$ws.Range("K77").Value = 34927.5
$ws.Range("L77").Value = 116942.004
$ws.Range("M77").Value = -30559.5
$ws.Range("N77").Value = -125678.004

$ws.Range("H106").Value = 69000
$ws.Range("J106").Value = 69000
$ws.Range("L106").Value = 69000
$ws.Range("N106").Value = -71524

$ws.Range("H107").Value = 352.04544
$ws.Range("I107").Value = 339.5625
$ws.Range("J107").Value = 385.33334
$ws.Range("K107").Value = 339.5625
$ws.Range("L107").Value = 385.33334
$ws.Range("M107").Value = 1580.4375
$ws.Range("N107").Value = -4225.33334

$ws.Range("H122").Value = 6394
$ws.Range("I122").Value = 2298
$ws.Range("K122").Value = 6894
$ws.Range("M122").Value = -4444

$ws.Range("H132").Value = 3836.2
$ws.Range("I132").Value = 3277.2727
$ws.Range("K132").Value = 9831.8181
$ws.Range("M132").Value = -7301.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 56.766666
$ws.Range("I2").Value = 35
$ws.Range("J2").Value = 58.32143
$ws.Range("K2").Value = 210
$ws.Range("L2").Value = 349.92858
$ws.Range("M2").Value = -97
$ws.Range("N2").Value = -575.92858

$ws.Range("H47").Value = 587.7778
$ws.Range("I47").Value = 458
$ws.Range("J47").Value = 750
$ws.Range("K47").Value = 1374
$ws.Range("L47").Value = 2250
$ws.Range("M47").Value = -943
$ws.Range("N47").Value = -3112

$ws.Range("H70").Value = 5151.5
$ws.Range("I70").Value = 924.6667
$ws.Range("J70").Value = 5996.8667
$ws.Range("K70").Value = 2774.0001
$ws.Range("L70").Value = 17990.6001
$ws.Range("M70").Value = -2459.0001
$ws.Range("N70").Value = -18620.6001

$ws.Range("H73").Value = 5151.5
$ws.Range("I73").Value = 924.6667
$ws.Range("J73").Value = 5996.8667
$ws.Range("K73").Value = 2774.0001
$ws.Range("L73").Value = 17990.6001
$ws.Range("M73").Value = -1682.0001
$ws.Range("N73").Value = -20174.6001

$ws.Range("H87").Value = 6674.75
$ws.Range("I87").Value = 2773.4285
$ws.Range("J87").Value = 8775.462
$ws.Range("K87").Value = 8320.2855
$ws.Range("L87").Value = 26326.386
$ws.Range("M87").Value = -7072.2855
$ws.Range("N87").Value = -28822.386

$ws.Range("H90").Value = 6674.75
$ws.Range("I90").Value = 2773.4285
$ws.Range("J90").Value = 8775.462
$ws.Range("K90").Value = 24960.8565
$ws.Range("L90").Value = 78979.158
$ws.Range("M90").Value = -18720.8565
$ws.Range("N90").Value = -91459.158

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5820.6978
$ws.Range("I70").Value = 5548.231
$ws.Range("K70").Value = 5548.231
$ws.Range("M70").Value = -5278.231

$ws.Range("H73").Value = 5820.6978
$ws.Range("I73").Value = 5548.231
$ws.Range("K73").Value = 5548.231
$ws.Range("M73").Value = -4612.231

$ws.Range("H101").Value = 41657
$ws.Range("J101").Value = 41657
$ws.Range("L101").Value = 41657
$ws.Range("N101").Value = -48147

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 400294.6
$ws.Range("J55").Value = 422
$ws.Range("L55").Value = 422
$ws.Range("N55").Value = -768

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 11451.5
$ws.Range("I55").Value = 2850
$ws.Range("J55").Value = 20053
$ws.Range("K55").Value = 2850
$ws.Range("L55").Value = 20053
$ws.Range("M55").Value = -2573
$ws.Range("N55").Value = -20607
